$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete rows (old rows 14-17, Resolving-Mac as target cluster)
$ws.Rows("14:17").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb1"
$ws.Range("C2").Value = "Ephb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.46510533333333
$ws.Range("H2").Value = 31.395316
$ws.Range("I2").Value = 0.5554075997074865
$ws.Range("J2").Value = 0.5554075997074865
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.168144
$ws.Range("N2").Value = 0.504432
$ws.Range("O2").Value = 0.05446245276675245
$ws.Range("P2").Value = 0.05446245276675245
$ws.Range("Q2").Value = 1.759644671168
$ws.Range("R2").Value = 15.836802040512
$ws.Range("S2").Value = 0.03024886016536434
$ws.Range("T2").Value = 0.03024886016536434

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb1"
$ws.Range("C3").Value = "Ephb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.46510533333333
$ws.Range("H3").Value = 31.395316
$ws.Range("I3").Value = 0.5554075997074865
$ws.Range("J3").Value = 0.5554075997074865
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.279928333333333
$ws.Range("N3").Value = 3.839785
$ws.Range("O3").Value = 0.4145734394268892
$ws.Range("P3").Value = 0.4145734394268892
$ws.Range("Q3").Value = 13.39458482745111
$ws.Range("R3").Value = 120.55126344706
$ws.Range("S3").Value = 0.2302572388945656
$ws.Range("T3").Value = 0.2302572388945656

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb1"
$ws.Range("C4").Value = "Ephb6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.46510533333333
$ws.Range("H4").Value = 31.395316
$ws.Range("I4").Value = 0.5554075997074865
$ws.Range("J4").Value = 0.5554075997074865
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.639265666666667
$ws.Range("N4").Value = 4.917797
$ws.Range("O4").Value = 0.5309641078063584
$ws.Range("P4").Value = 0.5309641078063583
$ws.Range("Q4").Value = 17.15508787098356
$ws.Range("R4").Value = 154.395790838852
$ws.Range("S4").Value = 0.2949015006475566
$ws.Range("T4").Value = 0.2949015006475565

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb1"
$ws.Range("C5").Value = "Ephb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.754308333333334
$ws.Range("H5").Value = 17.262925
$ws.Range("I5").Value = 0.3053945925621632
$ws.Range("J5").Value = 0.3053945925621632
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.168144
$ws.Range("N5").Value = 0.504432
$ws.Range("O5").Value = 0.05446245276675245
$ws.Range("P5").Value = 0.05446245276675245
$ws.Range("Q5").Value = 0.9675524204000001
$ws.Range("R5").Value = 8.707971783600001
$ws.Range("S5").Value = 0.01663253857263842
$ws.Range("T5").Value = 0.01663253857263842

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb1"
$ws.Range("C6").Value = "Ephb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.754308333333334
$ws.Range("H6").Value = 17.262925
$ws.Range("I6").Value = 0.3053945925621632
$ws.Range("J6").Value = 0.3053945925621632
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.279928333333333
$ws.Range("N6").Value = 3.839785
$ws.Range("O6").Value = 0.4145734394268892
$ws.Range("P6").Value = 0.4145734394268892
$ws.Range("Q6").Value = 7.365102274569446
$ws.Range("R6").Value = 66.28592047112501
$ws.Range("S6").Value = 0.1266084866208695
$ws.Range("T6").Value = 0.1266084866208695

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb1"
$ws.Range("C7").Value = "Ephb6"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.754308333333334
$ws.Range("H7").Value = 17.262925
$ws.Range("I7").Value = 0.3053945925621632
$ws.Range("J7").Value = 0.3053945925621632
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.639265666666667
$ws.Range("N7").Value = 4.917797
$ws.Range("O7").Value = 0.5309641078063584
$ws.Range("P7").Value = 0.5309641078063583
$ws.Range("Q7").Value = 9.432840086247223
$ws.Range("R7").Value = 84.89556077622501
$ws.Range("S7").Value = 0.1621535673686553
$ws.Range("T7").Value = 0.1621535673686553

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Efnb1"
$ws.Range("C8").Value = "Ephb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.146766
$ws.Range("H8").Value = 6.440298
$ws.Range("I8").Value = 0.1139338891693565
$ws.Range("J8").Value = 0.1139338891693565
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.168144
$ws.Range("N8").Value = 0.504432
$ws.Range("O8").Value = 0.05446245276675245
$ws.Range("P8").Value = 0.05446245276675245
$ws.Range("Q8").Value = 0.360965822304
$ws.Range("R8").Value = 3.248692400736
$ws.Range("S8").Value = 0.006205119057418489
$ws.Range("T8").Value = 0.006205119057418489

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Efnb1"
$ws.Range("C9").Value = "Ephb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.146766
$ws.Range("H9").Value = 6.440298
$ws.Range("I9").Value = 0.1139338891693565
$ws.Range("J9").Value = 0.1139338891693565
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.279928333333333
$ws.Range("N9").Value = 3.839785
$ws.Range("O9").Value = 0.4145734394268892
$ws.Range("P9").Value = 0.4145734394268892
$ws.Range("Q9").Value = 2.747706628436667
$ws.Range("R9").Value = 24.72935965593
$ws.Range("S9").Value = 0.04723396430022214
$ws.Range("T9").Value = 0.04723396430022214

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efnb1"
$ws.Range("C10").Value = "Ephb6"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.146766
$ws.Range("H10").Value = 6.440298
$ws.Range("I10").Value = 0.1139338891693565
$ws.Range("J10").Value = 0.1139338891693565
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.639265666666667
$ws.Range("N10").Value = 4.917797
$ws.Range("O10").Value = 0.5309641078063584
$ws.Range("P10").Value = 0.5309641078063583
$ws.Range("Q10").Value = 3.519119798167333
$ws.Range("R10").Value = 31.672078183506
$ws.Range("S10").Value = 0.06049480581171591
$ws.Range("T10").Value = 0.0604948058117159

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Efnb1"
$ws.Range("C11").Value = "Ephb6"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.4760280000000001
$ws.Range("H11").Value = 1.428084
$ws.Range("I11").Value = 0.02526391856099382
$ws.Range("J11").Value = 0.02526391856099382
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.168144
$ws.Range("N11").Value = 0.504432
$ws.Range("O11").Value = 0.05446245276675245
$ws.Range("P11").Value = 0.05446245276675245
$ws.Range("Q11").Value = 0.08004125203200001
$ws.Range("R11").Value = 0.720371268288
$ws.Range("S11").Value = 0.001375934971331207
$ws.Range("T11").Value = 0.001375934971331207

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Efnb1"
$ws.Range("C12").Value = "Ephb6"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.4760280000000001
$ws.Range("H12").Value = 1.428084
$ws.Range("I12").Value = 0.02526391856099382
$ws.Range("J12").Value = 0.02526391856099382
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.279928333333333
$ws.Range("N12").Value = 3.839785
$ws.Range("O12").Value = 0.4145734394268892
$ws.Range("P12").Value = 0.4145734394268892
$ws.Range("Q12").Value = 0.6092817246600001
$ws.Range("R12").Value = 5.48353552194
$ws.Range("S12").Value = 0.01047374961123203
$ws.Range("T12").Value = 0.01047374961123203

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Efnb1"
$ws.Range("C13").Value = "Ephb6"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4760280000000001
$ws.Range("H13").Value = 1.428084
$ws.Range("I13").Value = 0.02526391856099382
$ws.Range("J13").Value = 0.02526391856099382
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.639265666666667
$ws.Range("N13").Value = 4.917797
$ws.Range("O13").Value = 0.5309641078063584
$ws.Range("P13").Value = 0.5309641078063583
$ws.Range("Q13").Value = 0.7803363567720001
$ws.Range("R13").Value = 7.023027210948001
$ws.Range("S13").Value = 0.01341423397843058
$ws.Range("T13").Value = 0.01341423397843058
